$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ClientInfo")
$ws2 = $wb.Worksheets.Item("ClientAccountInfo")

# reference style (plain/default, same as the other numeric data cells)
$plainStyle = $ws1.Range("A2").Style

# New client record 1 (row 3 on ClientInfo)
$ws1.Range("A3").Value = 359
$ws1.Range("B3").Value = "'123"
$ws1.Range("B3").Style = $plainStyle
$ws1.Range("C3").Value = "'123"
$ws1.Range("C3").Style = $plainStyle
$ws1.Range("D3").Value = 123

# New client record 2 (row 4 on ClientInfo)
$ws1.Range("A4").Value = 160
$ws1.Range("B4").Value = "'1234"
$ws1.Range("B4").Style = $plainStyle
$ws1.Range("C4").Value = "'1234"
$ws1.Range("C4").Style = $plainStyle
$ws1.Range("D4").Value = 1234

# Corresponding login/password rows on ClientAccountInfo
$plainStyle2 = $ws2.Range("B2").Style

$ws2.Range("B3").Value = "'1"
$ws2.Range("B3").Style = $plainStyle2
$ws2.Range("C3").Value = "'1"
$ws2.Range("C3").Style = $plainStyle2

$ws2.Range("B4").Value = "'1234"
$ws2.Range("B4").Style = $plainStyle2
$ws2.Range("C4").Value = "'1234"
$ws2.Range("C4").Style = $plainStyle2
